# Update "想去人数" (attendee count) figures for two events that appear on
# both the "展览" sheet and the "全部类型" sheet (which mirrors all events).
#
#   展览 / 全部类型  row 2 (丽水·2024首届TCT国风动漫游戏嘉年华): F2  1446 -> 1448
#   展览 / 全部类型  row 4 (丽水·第三届HP国风动漫游戏嘉年华):   F4  14   -> 15

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 1448
    $ws.Range("F4").Value = 15
}
